$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52; this shifts the existing row 52 (and
# everything below it) down by one row, matching the target diff where
# every row from the old row 52 onward now lives one row lower.
$ws.Rows("52").Insert()

# Populate the newly inserted row 52 with the new data record.
$ws.Cells.Item(52, 1).Value2 = 11
$ws.Cells.Item(52, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(52, 3).Value2 = "Bíobío"
$ws.Cells.Item(52, 4).Value2 = 44791
$ws.Cells.Item(52, 5).Value2 = 8
$ws.Cells.Item(52, 6).Value2 = "Fruta"
$ws.Cells.Item(52, 7).Value2 = 100108
$ws.Cells.Item(52, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(52, 9).Value2 = 100108002
$ws.Cells.Item(52, 10).Value2 = "Mango"
$ws.Cells.Item(52, 11).Value2 = "Sin especificar"
$ws.Cells.Item(52, 12).Value2 = "Primera"
$ws.Cells.Item(52, 13).Value2 = 100
$ws.Cells.Item(52, 14).Value2 = 9500
$ws.Cells.Item(52, 15).Value2 = 10000
$ws.Cells.Item(52, 16).Value2 = 9750
$ws.Cells.Item(52, 17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(52, 18).Value2 = "Brasil"
$ws.Cells.Item(52, 19).Value2 = 2438
$ws.Cells.Item(52, 20).Value2 = 4
